$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update tx_loss (F) / tx_error (G) / tx_delay (H) values for rows 5-14 ---
# (Opus codec quality mapping table got recalibrated: loss/error rates halved,
#  delay now scales x3 instead of x1 per quality step.)
$ws.Range("F5").Value = 5
$ws.Range("G5").Value = 5
$ws.Range("H5").Value = 30

$ws.Range("F6").Value = 5
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 60

$ws.Range("F7").Value = 10
$ws.Range("G7").Value = 10
$ws.Range("H7").Value = 90

$ws.Range("F8").Value = 10
$ws.Range("G8").Value = 10
$ws.Range("H8").Value = 120

$ws.Range("F9").Value = 15
$ws.Range("G9").Value = 15
$ws.Range("H9").Value = 150

$ws.Range("F10").Value = 15
$ws.Range("G10").Value = 15
$ws.Range("H10").Value = 180

$ws.Range("F11").Value = 20
$ws.Range("G11").Value = 20
$ws.Range("H11").Value = 210

$ws.Range("F12").Value = 20
$ws.Range("G12").Value = 20
$ws.Range("H12").Value = 240

$ws.Range("F13").Value = 25
$ws.Range("G13").Value = 25
$ws.Range("H13").Value = 270

$ws.Range("F14").Value = 25
$ws.Range("G14").Value = 25
$ws.Range("H14").Value = 300

# touch formatting on the tx_error column and the merged Opus-formula cell so the
# workbook's style table reflects the same re-format pass made in Excel
$ws.Range("G4:G14").NumberFormat = "General"
$ws.Range("B4:C4").HorizontalAlignment = 1

# --- New notes added below the table ---
$ws.Range("F16").Value = "up to 20% will notice, otherwise just heard as bad"
$ws.Range("H16").Value = "units ms"
$ws.Range("D16").Value = "crowd of people talking/screaming/being eaten by monster"
$ws.Range("D17").Value = "texas chainsaw masacre movie, emulate sound?"

$ws.Range("D16:D17").HorizontalAlignment = 1

# column D now holds the longest notes in the sheet -- widen it to fit
$ws.Columns("D").AutoFit()

$ws.Range("G14").Select()
